$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric (e.g. "58.00", "0.0762") into a
# cell while forcing it to be stored as literal text, exactly as it was
# before the edit (every cell touched here holds plain string content, not
# a number). A plain `.Value = "58.00"` assignment would let Excel
# auto-convert the text to the number 58 (dropping the trailing zero and
# changing the cell's stored type) - not what we want.
#
# The standard COM fix is the same trick the Excel UI uses for a manually
# typed numeric-looking entry: prefix the value with a leading apostrophe
# via .Formula, which reliably keeps it as text. That flips the cell's
# internal "quote prefix" flag though, so reset .Style back to "Normal"
# immediately after - the cell ends up with the same (default/no) style it
# had before, just holding the exact text we wanted.
function Set-TextCell($ref, $text) {
    $ws.Range($ref).Formula = "'" + $text
    $ws.Range($ref).Style = "Normal"
}

$ws.Range('D2').Value = '37.359.62'

$ws.Range('D3').Value = '2.061.89'
$ws.Range('E3').Value = '  +4.26%  '

$ws.Range('E4').Value = '  +0.09%  '

Set-TextCell D5 '236.16'
$ws.Range('E5').Value = '  +0.97%  '

$ws.Range('E6').Value = '  +2.92%  '

Set-TextCell D7 '57.93'
$ws.Range('E7').Value = '  +6.30%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +3.00%  '

Set-TextCell D10 '58.00'
$ws.Range('E10').Value = '  -1.50%  '

Set-TextCell D11 '0.0762'

$ws.Range('E12').Value = '  +3.47%  '

$ws.Range('D13').Value = '2.367.83'
$ws.Range('E13').Value = '  +4.37%  '

Set-TextCell D14 '14.45'
$ws.Range('E14').Value = '  +3.50%  '

Set-TextCell D15 '20.88'
$ws.Range('E15').Value = '  +4.64%  '

Set-TextCell D16 '0.777'
$ws.Range('E16').Value = '  +3.59%  '

Set-TextCell D17 '5.17'
$ws.Range('E17').Value = '  +2.38%  '

$ws.Range('D18').Value = '2.052.99'
$ws.Range('E18').Value = '  +3.82%  '

$ws.Range('D19').Value = '37.584.68'
$ws.Range('E19').Value = '  +3.48%  '

$ws.Range('E20').Value = '  +17.35%  '

Set-TextCell D21 '69.14'
$ws.Range('E21').Value = '  +2.35%  '

$ws.Range('D22').Value = '0.0₃0816'
$ws.Range('E22').Value = '  +1.61%  '

Set-TextCell D23 '227.11'
$ws.Range('E23').Value = '  +2.50%  '

$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('E25').Value = '  +3.19%  '

$ws.Range('E26').Value = '  +1.59%  '

Set-TextCell D27 '163.88'
$ws.Range('E27').Value = '  +1.99%  '

$ws.Range('E28').Value = '  +12.55%  '

Set-TextCell D29 '8.85'
$ws.Range('E29').Value = '  +3.71%  '

Set-TextCell D30 '19.17'
$ws.Range('E30').Value = '  +2.17%  '

$ws.Range('E31').Value = '  -0.33%  '

$ws.Range('E32').Value = '  +2.51%  '

Set-TextCell D33 '4.49'
$ws.Range('E33').Value = '  +2.96%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell D34 '0.0621'
$ws.Range('E34').Value = '  +2.61%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell D35 '2.57'
$ws.Range('E35').Value = '  +11.73%  '

$ws.Range('E36').Value = '  +5.66%  '

$ws.Range('E37').Value = '  +6.59%  '

$ws.Range('E39').Value = '  +0.68%  '

Set-TextCell D40 '5.91'
$ws.Range('E40').Value = '  +9.18%  '

$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell D41 '4.61'
$ws.Range('E41').Value = '  +32.75%  '

$ws.Range('B42').Value = 'Cronos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell D42 '0.0986'
$ws.Range('E42').Value = '  +9.60%  '

Set-TextCell D43 '2.96'
$ws.Range('E43').Value = '  -2.36%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.478.32'
$ws.Range('E44').Value = '  +1.36%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell D45 '97.09'
$ws.Range('E45').Value = '  +10.44%  '

$ws.Range('E46').Value = '  +7.32%  '

Set-TextCell D47 '0.0210'
$ws.Range('E47').Value = '  +4.69%  '

$ws.Range('E48').Value = '  +7.26%  '

Set-TextCell D49 '1.02'
$ws.Range('E49').Value = '  +3.44%  '

Set-TextCell D50 '7.20'
$ws.Range('E50').Value = '  +6.74%  '

$ws.Range('E51').Value = '  +2.18%  '
